$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared "Zone1" string to "Zone1_1" (affects row 2, F2)
$ws.Range("F2").Value = "Zone1_1"

# Row 3 (HostName2): ZoneID/RackID default now reuse Zone1_1 / Rack1
$ws.Range("F3").Value = "Zone1_1"
$ws.Range("G3").Value = "Rack1"

# Row 4 (HostName3): ZoneID/RackID default now reuse Zone1_1 / Rack1
$ws.Range("F4").Value = "Zone1_1"
$ws.Range("G4").Value = "Rack1"

# Make the "Component Purpose" header run bold (matches sibling header styling)
$ws.Range("O1").Characters(1, 17).Font.Bold = $true
$ws.Range("O1").Characters(1, 17).Font.Size = 12
$ws.Range("O1").Characters(1, 17).Font.Name = "宋体"

# Move the active selection to A2
$ws.Range("A2").Select() | Out-Null
